$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Work bottom-up (paragraph 3, then 2, then 1) so that edits which add
# paragraphs never shift the index of a paragraph we haven't processed
# yet.
# ---------------------------------------------------------------------

# --- Original paragraph 3 (bookmark + "Un handyWorker...") becomes the
# "Permitir en el pattern..." bullet; three brand new bullets (Dejar los
# atributos / Un mensaje solo / Un handyWorker... with the bookmark moved
# to the end) are appended right after it. ---
$p3 = $d.Paragraphs(3)
$xmlT7 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Permitir en el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pattern</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> todos</w:t></w:r><w:r><w:t xml:space="preserve"> los</w:t></w:r><w:r><w:t xml:space="preserve"> formato</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> posibles de los emails. En los servicios se restringen los patrones exclusivos de administrador</w:t></w:r></w:p>'
$p3.Range.InsertXML($xmlT7)

$p3 = $d.Paragraphs(3)
$endOfP3 = $p3.Range.End
$insertPoint = $d.Range($endOfP3, $endOfP3)
$xmlTail = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Dejar los atributos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>warranty</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>category</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> como tipo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>string</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, ya que podemos entender que el buscador sigue la filosofía de un buscador puro y duro. Si cambia el nombre de una categoría no hay ningún problema, ya que en un buscador podemos tener parámetros de búsqueda que no encuentren nada.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Un mensaje solo se elimina de la base de datos cuando no esté referenciado en ninguna carpeta de ningún actor.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">¿Un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>handyWorker</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> podría solicitar más de una vez una tarea? NO HAY CLASE ASOCIACIÓN. Podría interpretarse como entidad ya que puede cambiar el estado a REJECTED y volverlo a solicitar, cambiando así el estado y por lo tanto en ese caso sería una entidad. Pero esto depende de la interpretación de los requisitos.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$insertPoint.InsertXML($xmlTail)

# --- Original paragraph 2 ("Para un handyWorker...") just gets its runs
# swapped for the new "En la entrega de servicios..." sentence (via
# InsertXML, since Range.Text only overwrites the first run and would
# leave the old proofErr/run tail behind) -- this keeps the paragraph
# count (and therefore every other index) unchanged. ---
$p2 = $d.Paragraphs(2)
$xmlT6 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>En la entrega de servicios tratar de solucionar el siguiente problema con las categorías. En el caso de que se añada un nuevo idioma al sistema, inicialmente todas las categorías tienen ese nuevo idioma a nulo, lo que puede dar error.</w:t></w:r></w:p>'
$p2.Range.InsertXML($xmlT6)

# --- Original paragraph 1 ("Un mensaje puede estar...") turns into an
# empty ind-only paragraph, followed by the "Alternativas..." bullet and
# its three ilvl=1 sub-bullets. Processed last since it shifts every
# later paragraph index. ---
$p1 = $d.Paragraphs(1)
$xmlBatch1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Alternativas posibles de las asociaciones </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Endorsement</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Endorsable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Dejar las navegabilidades como están</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Poner las navegabilidades bidireccionales</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Poner las navegabilidades en sentido contrario</w:t></w:r></w:p>'
$p1.Range.InsertXML($xmlBatch1)

